$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = -7
$ws.Range("F4").Value = 5
$ws.Range("F7").Value = -11
$ws.Range("F8").Value = -5
$ws.Range("F9").Value = 5
$ws.Range("F15").Value = -3
$ws.Range("F16").Value = -1
$ws.Range("F19").Value = -4
$ws.Range("F20").Value = 3
$ws.Range("F24").Value = -5
$ws.Range("F25").Value = -3
$ws.Range("F27").Value = 2
$ws.Range("F35").Value = 6
$ws.Range("F37").Value = -2
$ws.Range("F39").Value = -7
$ws.Range("F43").Value = -7
$ws.Range("F45").Value = -4
$ws.Range("F47").Value = -3
$ws.Range("F48").Value = 3
$ws.Range("F55").Value = 3
$ws.Range("F56").Value = 1
$ws.Range("F58").Value = -3
